$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80-88 down to 81-89
$ws.Rows.Item(80).Insert()

# Populate the new row 80 with the new data record, following the same
# pattern as the surrounding rows for the columns that stay constant.
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "Vega Modelo de Temuco"
$ws.Range("C80").Value = "La Araucanía"
$ws.Range("D80").Value = 44491
$ws.Range("E80").Value = 9
$ws.Range("F80").Value = 100112012
$ws.Range("G80").Value = "Espinaca"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 65
$ws.Range("K80").Value = 8000
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = 8000
$ws.Range("N80").Value = "$/docena de atados"
$ws.Range("O80").Value = "Región de La Araucanía"
$ws.Range("P80").Value = 2667
$ws.Range("Q80").Value = 3
$ws.Range("R80").Value = "Hortaliza"
